$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New observation row appended at the bottom of the sheet (row 9)

$ws.Range("A9").Value = 112126339
$ws.Range("B9").Value = 90332
$ws.Range("C9").Value = "Ovaliderad"
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 4769
$ws.Range("F9").Value = "Svavelriska"
$ws.Range("G9").Value = "Lactarius scrobiculatus"
$ws.Range("H9").Value = "(Scop.:Fr.) Fr."

# Present-but-empty text cell (matches the blank <is><t/></is> placeholder
# pattern already used elsewhere in the sheet, e.g. column I on other rows)
$ws.Range("I9").Value = "'"

$ws.Range("P9").Value = "Brännberget (N om reservatet), Vb"
$ws.Range("Q9").Value = 760438.3549515785
$ws.Range("R9").Value = 7210190.261124903
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = "Västerbotten"
$ws.Range("U9").Value = "Skellefteå"
$ws.Range("V9").Value = "Västerbotten"
$ws.Range("W9").Value = "Skellefteå socken"

# Dates/times are stored as plain text in this sheet, not native Excel
# dates, so force text formatting before assigning to avoid auto
# date-serial conversion.
$ws.Range("Y9:AB9").NumberFormat = "@"
$ws.Range("Y9").Value = "2023-09-15"
$ws.Range("Z9").Value = "00:00"
$ws.Range("AA9").Value = "2023-09-15"
$ws.Range("AB9").Value = "00:00"

$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false

$ws.Range("AT9").Value = "'"

$ws.Range("AW9").Value = "Emil Larsson"
$ws.Range("AX9").Value = "Emil Larsson, Carl Jansson"

$ws.Range("AY9").Value = "'"
